# Updated DPM integration testfixture with hierarchy node labels
#
# The fixture's hierarchy-node identifiers (the "ID" column on the
# CodeSchemes sheet, and the "ID"/"SUBCODESCHEME" columns on the Codes
# sheet, plus the "ID" column on the Extensions sheet) were regenerated.
# Column A on each sheet is sized to fit its (longest) ID, so after the
# values change we re-apply the column's best-fit width.

$wb = $excel.ActiveWorkbook

# --- Sheet "CodeSchemes" ------------------------------------------------
$ws1 = $wb.Worksheets.Item("CodeSchemes")
$ws1.Range("A2").Value = "9560e39a-7e35-46ea-b8af-afbb30a91a4b"
$ws1.Columns.Item(1).ColumnWidth = 29

# --- Sheet "Codes" --------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Codes")
$ws2.Range("A2").Value = "9c159c69-307c-42f7-a146-9b8d940b064f"
$ws2.Range("K2").Value = "a68e6334-aedf-4a28-9964-057446f62d8c"

$ws2.Range("A3").Value = "750f8958-8106-4f72-bff6-62ec413eeb7c"
$ws2.Range("K3").Value = "f01bc58e-d294-4d60-a62c-ef0670327a6c"

$ws2.Range("A4").Value = "27cb809b-3fa5-487a-bacc-db636bddda1a"
$ws2.Range("K4").Value = "b2da0f1c-f28c-4559-a5f8-72b0379930d5"

$ws2.Range("A5").Value = "4959b7dc-56d3-4d29-80bf-c4c4788d42ce"
$ws2.Range("K5").Value = "54afa059-9b2e-4ae8-b566-f001a23f4e21"

$ws2.Range("A6").Value = "b2b7ae8c-3d91-4e8c-84a2-cccd538874fb"
$ws2.Range("K6").Value = "5e4d6018-a396-4bb9-b75d-9fc0f91d4d50"

$ws2.Range("A7").Value = "a3b3d2ee-1a83-4e12-9edb-819c2cbcbc16"
$ws2.Range("K7").Value = "ac63e4ba-8c75-4a63-a4e5-d4149d960bf2"

$ws2.Range("A8").Value = "d3628223-4b6a-4055-982c-32b04c45f38b"
$ws2.Range("K8").Value = "5c8de2aa-8d6b-4814-937f-920da3190315"

$ws2.Columns.Item(1).ColumnWidth = 33.42857142857143

# --- Sheet "Extensions" ----------------------------------------------------
$ws3 = $wb.Worksheets.Item("Extensions")
$ws3.Range("A2").Value = "854696e4-fb38-4867-a986-65dd00ead244"
$ws3.Columns.Item(1).ColumnWidth = 34.42857142857143
